# Added Filter for Customer in the Specific Pricing Feature
#
# Semantic changes applied:
#  1. The "Web Data 12" label (shared across the header row of every sheet)
#     is renamed to "Web Data 3".
#  2. Selection/active-cell bookkeeping is refreshed: the Pricing and Order
#     sheets move their selection to C1 (the newly renamed filter label),
#     and the workbook's active sheet moves from "Order" back to "Pricing".

$wb = $excel.ActiveWorkbook

# 1) Rename the shared "Web Data 12" -> "Web Data 3" label in C1 on every sheet.
foreach ($name in @("Pricing", "Product", "Customer", "Order")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value = "Web Data 3"
}

# 2) Refresh the Order sheet's selection (was H3 -> now C1), matching the
#    new filter-label cell. This also activates Order momentarily.
$wsOrder = $wb.Worksheets.Item("Order")
$wsOrder.Range("C1").Select()

# 3) Finally, move back to the Pricing sheet (the workbook's active tab),
#    and refresh its selection (was H1 -> now C1). Doing this last makes
#    Pricing the active/selected sheet on save, matching the target state.
$wsPricing = $wb.Worksheets.Item("Pricing")
$wsPricing.Activate()
$wsPricing.Range("C1").Select()
